# Resolve timestamps when import file
# Rewrites the "Waktu" (H) column values with fresh import timestamps,
# renumbers the "ID" (B) column for the three batches of 5 guests each,
# and appends two more full batches of 5 guest rows (rows 10-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fixed values shared by every data row (unchanged across the whole table)
$colC = "sdwdwd"
$colD = "wdwdw"
$colE = "dwdwd"
$colF = "wdwdwdw"
$phone = 4242424242

# Three import batches (newest first), each with its own "Waktu" stamp.
# Batch rows are written oldest-appearing-last, matching the worksheet order:
#   rows 2-6   -> newest batch, IDs 11-15, "Friday, 06-12-2024 01:18:09"
#   rows 7-11  -> middle batch, IDs 6-10,  "Friday, 06-12-2024 01:14:47"
#   rows 12-16 -> oldest batch, IDs 1-5,   "Friday, 06-12-2024 01:12:40"
$batches = @(
    @{ Waktu = "Friday, 06-12-2024 01:18:09"; Ids = @(11, 12, 13, 14, 15) },
    @{ Waktu = "Friday, 06-12-2024 01:14:47"; Ids = @(6, 7, 8, 9, 10) },
    @{ Waktu = "Friday, 06-12-2024 01:12:40"; Ids = @(1, 2, 3, 4, 5) }
)

# Template row used to stamp the existing data-row formatting (style index 1)
# onto the freshly appended rows (10-16) so they look like the rest of the table.
$templateRow = $ws.Range("A9:H9")
$lastExistingRow = 9

$row = 2
$no = 1
foreach ($batch in $batches) {
    foreach ($id in $batch.Ids) {
        if ($row -gt $lastExistingRow) {
            $templateRow.Copy()
            $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 8)).PasteSpecial(-4122)
        }

        $ws.Cells.Item($row, 1).Value = $no
        $ws.Cells.Item($row, 2).Value = $id
        $ws.Cells.Item($row, 3).Value = $colC
        $ws.Cells.Item($row, 4).Value = $colD
        $ws.Cells.Item($row, 5).Value = $colE
        $ws.Cells.Item($row, 6).Value = $colF
        $ws.Cells.Item($row, 7).Value = $phone
        $ws.Cells.Item($row, 8).Value = $batch.Waktu

        $no = $no + 1
        $row = $row + 1
    }
}
